$wb = $excel.ActiveWorkbook

# --- Rename the "Weight" label to "BMI" (row label + column header) ---
$rho = $wb.Worksheets.Item("Rho")
$pval = $wb.Worksheets.Item("Pval")

$rho.Range("A3").Value = "BMI"
$rho.Range("C1").Value = "BMI"
$pval.Range("A3").Value = "BMI"
$pval.Range("C1").Value = "BMI"

# --- Rho sheet: widen columns C and D by one unit (14.71... -> 15.71...) ---
$rho.Columns.Item(3).ColumnWidth = 14.877604166666666
$rho.Columns.Item(4).ColumnWidth = 14.877604166666666

# --- Rho sheet: updated correlation coefficients ---
$rho.Range("C2").Value = -0.046295178066277004
$rho.Range("D2").Value = -0.4814930477872299
$rho.Range("E2").Value = -0.25485144331159998

$rho.Range("B3").Value = -0.046295178066277004
$rho.Range("D3").Value = -0.006878549834324503
$rho.Range("E3").Value = 0.013940261464136558

$rho.Range("B4").Value = -0.4814930477872299
$rho.Range("C4").Value = -0.006878549834324503
$rho.Range("E4").Value = -0.25619521300354525

$rho.Range("B5").Value = -0.25485144331159998
$rho.Range("C5").Value = 0.013940261464136558
$rho.Range("D5").Value = -0.25619521300354525

# --- Pval sheet: columns B, C, D narrowed by one unit ---
$pval.Columns.Item(2).ColumnWidth = 12.877604166666666
$pval.Columns.Item(3).ColumnWidth = 10.877604166666666
$pval.Columns.Item(4).ColumnWidth = 12.877604166666666

# --- Pval sheet: updated p-values ---
$pval.Range("C2").Value = 0.85995285403144428
$pval.Range("D2").Value = 0.050356403552665012
$pval.Range("E2").Value = 0.3235592506313606

$pval.Range("B3").Value = 0.85995285403144428
$pval.Range("D3").Value = 0.97909716700286242
$pval.Range("E3").Value = 0.95765122396067226

$pval.Range("B4").Value = 0.050356403552665012
$pval.Range("C4").Value = 0.97909716700286242
$pval.Range("E4").Value = 0.32092613009253734

$pval.Range("B5").Value = 0.3235592506313606
$pval.Range("C5").Value = 0.95765122396067226
$pval.Range("D5").Value = 0.32092613009253734
